$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = "AlexNet Implementation(구현) by PyTorch"
$ws.Range("E4").Value = "https://teddylee777.github.io/pytorch/12-alexnet-implementation"

# Row 26
$ws.Range("D26").Value = "인공지능 음성 생성 연구: 음성 분류 솔루션"

# Row 32
$ws.Range("D32").Value = "Feature Selection VS Feature Extraction"
$ws.Range("E32").Value = "https://dodonam.tistory.com/387"

# Row 50
$ws.Range("D50").Value = "가장 쉬운 Monte Carlo dropout [추론의 정밀도 산정에 이용]"
$ws.Range("E50").Value = "http://incredible.egloos.com/7547375"

# Row 51
$ws.Range("D51").Value = "[html] 웹페이지에 이모지 넣는 방법"
$ws.Range("E51").Value = "https://bskyvision.com/entry/html-%EC%9B%B9%ED%8E%98%EC%9D%B4%EC%A7%80%EC%97%90-%EC%9D%B4%EB%AA%A8%EC%A7%80-%EB%84%A3%EB%8A%94-%EB%B0%A9%EB%B2%95"
